$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.850.26"
$ws.Range("E2").Value = "  -0.41%  "
$ws.Range("D3").Value = "1.629.87"
$ws.Range("E3").Value = "  -0.27%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.27"
$ws.Range("E5").Value = "  -0.39%  "
$ws.Range("E6").Value = "  -0.88%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.24"
$ws.Range("E8").Value = "  -0.90%  "
$ws.Range("E9").Value = "  -0.40%  "
$ws.Range("E10").Value = "  -1.10%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0880"
$ws.Range("E11").Value = "  -0.12%  "
$ws.Range("D12").Value = "1.860.69"
$ws.Range("D13").Value = "1.628.24"
$ws.Range("E13").Value = "  -0.29%  "
$ws.Range("E14").Value = "  -1.17%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.556"
$ws.Range("E15").Value = "  -1.29%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.95"
$ws.Range("E16").Value = "  -1.24%  "
$ws.Range("D17").Value = "27.865.50"
$ws.Range("E17").Value = "  -0.36%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "228.67"
$ws.Range("E18").Value = "  -1.50%  "
$ws.Range("E19").Value = "  +0.02%  "
$ws.Range("D20").Value = "0.0₃0719"
$ws.Range("E20").Value = "  -1.00%  "
$ws.Range("E21").Value = "  -0.25%  "
$ws.Range("E22").Value = "  -0.26%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.95"
$ws.Range("E23").Value = "  -4.56%  "
$ws.Range("E24").Value = "  -0.66%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "155.11"
$ws.Range("E25").Value = "  -0.24%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.92"
$ws.Range("E26").Value = "  -0.36%  "
$ws.Range("E27").Value = "  -0.27%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.46"
$ws.Range("E28").Value = "  -1.28%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.996"
$ws.Range("E29").Value = "  -0.30%  "
$ws.Range("E30").Value = "  -0.28%  "
$ws.Range("E31").Value = "  -0.20%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.40"
$ws.Range("E32").Value = "  -0.02%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.11"
$ws.Range("E33").Value = "  +0.12%  "
$ws.Range("D34").Value = "1.412.70"
$ws.Range("E34").Value = "  +0.35%  "
$ws.Range("E35").Value = "  +2.59%  "
$ws.Range("E36").Value = "  -2.77%  "
$ws.Range("E37").Value = "  -1.46%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0169"
$ws.Range("E38").Value = "  -0.93%  "
$ws.Range("E39").Value = "  -0.59%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.853"
$ws.Range("E40").Value = "  -1.62%  "
$ws.Range("E41").Value = "  -2.16%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "65.72"
$ws.Range("E42").Value = "  -1.68%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.82"
$ws.Range("E43").Value = "  -0.15%  "
$ws.Range("E44").Value = "  -0.58%  "
$ws.Range("D45").Value = "1.769.82"
$ws.Range("E45").Value = "  -0.34%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.11"
$ws.Range("E46").Value = "  -3.88%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "88.52"
$ws.Range("E47").Value = "  +0.40%  "
$ws.Range("E48").Value = "  +1.05%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.61"
$ws.Range("E50").Value = "  +1.47%  "
$ws.Range("E51").Value = "  -0.40%  "
